# Bug fix: connector line shapes were saved with an effectively
# zero line width (1 EMU, i.e. a:ln w="1"). Restore the intended
# hairline width of 1pt (12700 EMU) on every connector shape.
$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $shp = $s.Shapes.Item($i)
        if ($shp.Type -eq 9) {
            $shp.Line.Weight = 1
        }
    }
}
